$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 23.6
$ws.Range("I5").Value = 29.25
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 29.25
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 85.75
$ws.Range("N5").Value = -231
$ws.Range("H33").Value = 909.2759
$ws.Range("I33").Value = 352.76
$ws.Range("K33").Value = 352.76
$ws.Range("M33").Value = -123.76
$ws.Range("H64").Value = 41837.81
$ws.Range("J64").Value = 3549.2273
$ws.Range("L64").Value = 3549.2273
$ws.Range("N64").Value = -4045.2273
$ws.Range("H67").Value = 41837.81
$ws.Range("J67").Value = 3549.2273
$ws.Range("L67").Value = 3549.2273
$ws.Range("N67").Value = -5265.2273
$ws.Range("H111").Value = 6676290
$ws.Range("I111").Value = 18682.834
$ws.Range("J111").Value = 11114695
$ws.Range("K111").Value = 56048.50199999999
$ws.Range("L111").Value = 33344085
$ws.Range("M111").Value = -52981.50199999999
$ws.Range("N111").Value = -33350219
$ws.Range("H115").Value = 20354
$ws.Range("J115").Value = 570
$ws.Range("L115").Value = 1710
$ws.Range("N115").Value = -4844
$ws.Range("H129").Value = 4052.1292
$ws.Range("J129").Value = 1157
$ws.Range("L129").Value = 3471
$ws.Range("N129").Value = -13471
$ws.Range("H138").Value = 1994.9888
$ws.Range("I138").Value = 935.3137
$ws.Range("J138").Value = 3417.1843
$ws.Range("K138").Value = 2805.9411
$ws.Range("L138").Value = 10251.5529
$ws.Range("M138").Value = 2334.0589
$ws.Range("N138").Value = -20531.5529

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 127384.875
$ws.Range("I2").Value = 2915.8
$ws.Range("J2").Value = 334833.34
$ws.Range("K2").Value = 2915.8
$ws.Range("L2").Value = 334833.34
$ws.Range("M2").Value = -2802.8
$ws.Range("N2").Value = -335059.34
$ws.Range("H32").Value = 59614.934
$ws.Range("I32").Value = 11445.3
$ws.Range("J32").Value = 155954.2
$ws.Range("K32").Value = 11445.3
$ws.Range("L32").Value = 155954.2
$ws.Range("M32").Value = -11158.3
$ws.Range("N32").Value = -156528.2
$ws.Range("H45").Value = 103796
$ws.Range("I45").Value = 169696.67
$ws.Range("J45").Value = 4945
$ws.Range("K45").Value = 169696.67
$ws.Range("L45").Value = 4945
$ws.Range("M45").Value = -169319.67
$ws.Range("N45").Value = -5699
$ws.Range("H61").Value = 2102.2
$ws.Range("I61").Value = 2102.2
$ws.Range("K61").Value = 2102.2
$ws.Range("M61").Value = -1890.2
$ws.Range("H74").Value = 1345.7142
$ws.Range("I74").Value = 748.8
$ws.Range("J74").Value = 2141.6
$ws.Range("K74").Value = 748.8
$ws.Range("L74").Value = 2141.6
$ws.Range("M74").Value = 125.2
$ws.Range("N74").Value = -3889.6
$ws.Range("H77").Value = 1345.7142
$ws.Range("I77").Value = 748.8
$ws.Range("J77").Value = 2141.6
$ws.Range("K77").Value = 3744
$ws.Range("L77").Value = 10708
$ws.Range("M77").Value = 624
$ws.Range("N77").Value = -19444
$ws.Range("H92").Value = 20000
$ws.Range("J92").Value = 20000
$ws.Range("L92").Value = 20000
$ws.Range("N92").Value = -24992
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H116").Value = 127384.875
$ws.Range("I116").Value = 2915.8
$ws.Range("J116").Value = 334833.34
$ws.Range("K116").Value = 2915.8
$ws.Range("L116").Value = 334833.34
$ws.Range("M116").Value = -621.8000000000002
$ws.Range("N116").Value = -339421.34
$ws.Range("H136").Value = 2102.2
$ws.Range("I136").Value = 2102.2
$ws.Range("K136").Value = 6306.599999999999
$ws.Range("M136").Value = -3756.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 127384.875
$ws.Range("I3").Value = 2915.8
$ws.Range("J3").Value = 334833.34
$ws.Range("K3").Value = 2915.8
$ws.Range("L3").Value = 334833.34
$ws.Range("M3").Value = -2801.8
$ws.Range("N3").Value = -335061.34
$ws.Range("H108").Value = 34892
$ws.Range("J108").Value = 34892
$ws.Range("L108").Value = 34892
$ws.Range("N108").Value = -42572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22128.578
$ws.Range("I31").Value = 963.55554
$ws.Range("K31").Value = 963.55554
$ws.Range("M31").Value = -668.55554
$ws.Range("H34").Value = 22128.578
$ws.Range("I34").Value = 963.55554
$ws.Range("K34").Value = 963.55554
$ws.Range("M34").Value = -761.55554
$ws.Range("H63").Value = 42995
$ws.Range("J63").Value = 42995
$ws.Range("L63").Value = 42995
$ws.Range("N63").Value = -44367
$ws.Range("H66").Value = 42995
$ws.Range("J66").Value = 42995
$ws.Range("L66").Value = 128985
$ws.Range("N66").Value = -135849
$ws.Range("H81").Value = 39250
$ws.Range("J81").Value = 39250
$ws.Range("L81").Value = 39250
$ws.Range("N81").Value = -41246
$ws.Range("H84").Value = 39250
$ws.Range("J84").Value = 39250
$ws.Range("L84").Value = 117750
$ws.Range("N84").Value = -127734

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 288
$ws.Range("J34").Value = 900
$ws.Range("L34").Value = 2700
$ws.Range("N34").Value = -2868
$ws.Range("H131").Value = 818.5
$ws.Range("J131").Value = 852.06525
$ws.Range("L131").Value = 2556.19575
$ws.Range("N131").Value = -12636.19575

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 9196
$ws.Range("J123").Value = 9196
$ws.Range("L123").Value = 9196
$ws.Range("N123").Value = -14096

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2549.7144
$ws.Range("I7").Value = 1652.7858
$ws.Range("K7").Value = 1652.7858
$ws.Range("M7").Value = -1540.7858
$ws.Range("H40").Value = 64714.5
$ws.Range("I40").Value = 251550
$ws.Range("J40").Value = 2436
$ws.Range("K40").Value = 251550
$ws.Range("L40").Value = 2436
$ws.Range("M40").Value = -251414
$ws.Range("N40").Value = -2708
$ws.Range("H46").Value = 562843.4399999999
$ws.Range("I46").Value = 252
$ws.Range("J46").Value = 1266082.8
$ws.Range("K46").Value = 252
$ws.Range("L46").Value = 1266082.8
$ws.Range("M46").Value = -64
$ws.Range("N46").Value = -1266458.8
$ws.Range("H122").Value = 1942.5714
$ws.Range("I122").Value = 2066.6667
$ws.Range("J122").Value = 1849.5
$ws.Range("K122").Value = 6200.000100000001
$ws.Range("L122").Value = 5548.5
$ws.Range("M122").Value = -3750.000100000001
$ws.Range("N122").Value = -10448.5
$ws.Range("H126").Value = 2549.7144
$ws.Range("I126").Value = 1652.7858
$ws.Range("K126").Value = 4958.357400000001
$ws.Range("M126").Value = -2488.357400000001
$ws.Range("H134").Value = 65429
$ws.Range("J134").Value = 65429
$ws.Range("L134").Value = 65429
$ws.Range("N134").Value = -75569
$ws.Range("H136").Value = 2187.4
$ws.Range("J136").Value = 2132.8333
$ws.Range("L136").Value = 6398.499899999999
$ws.Range("N136").Value = -11498.4999
$ws.Range("H140").Value = 73445
$ws.Range("J140").Value = 73445
$ws.Range("L140").Value = 73445
$ws.Range("N140").Value = -83805

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 37664
$ws.Range("J16").Value = 37664
$ws.Range("L16").Value = 37664
$ws.Range("N16").Value = -38248
$ws.Range("H39").Value = 10960.5
$ws.Range("I39").Value = 5521
$ws.Range("K39").Value = 5521
$ws.Range("M39").Value = -5108
$ws.Range("H43").Value = 100000
$ws.Range("I43").Value = 100000
$ws.Range("K43").Value = 100000
$ws.Range("M43").Value = -99851
$ws.Range("H107").Value = 333866.34
$ws.Range("J107").Value = 500399.5
$ws.Range("L107").Value = 1501198.5
$ws.Range("N107").Value = -1505038.5
